$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 0.8136575222015381
$ws.Range("E2").Value = 5251.587009206085
$ws.Range("G2").Value = 0.1768077525252375
$ws.Range("H2").Value = 0.1562931651970017
$ws.Range("I2").Value = 0.1290241462329322
$ws.Range("J2").Value = 0.1222519465937708
$ws.Range("K2").Value = 0.1123494783187909
$ws.Range("L2").Value = 0.110599287216326
$ws.Range("M2").Value = 0.110599287216326
$ws.Range("N2").Value = 0.110599287216326
$ws.Range("O2").Value = 0.110599287216326
$ws.Range("P2").Value = 0.110599287216326
$ws.Range("Q2").Value = 0.110599287216326
$ws.Range("R2").Value = 0.110599287216326
$ws.Range("S2").Value = 0.110599287216326
$ws.Range("T2").Value = 0.110599287216326
$ws.Range("U2").Value = 0.110599287216326
$ws.Range("V2").Value = 0.1104908180339059
$ws.Range("W2").Value = 0.1103991267324654
$ws.Range("X2").Value = 0.1103990078872424
$ws.Range("Y2").Value = 0.1103701171385202
# Row 3
$ws.Range("C3").Value = 0.9103384017944336
$ws.Range("E3").Value = 5239.93958471834
$ws.Range("H3").Value = 0.161455175852402
$ws.Range("I3").Value = 0.1486736310454331
$ws.Range("J3").Value = 0.1348145567558006
$ws.Range("K3").Value = 0.1226590852641088
$ws.Range("L3").Value = 0.1113589360804189
$ws.Range("M3").Value = 0.1113589360804189
$ws.Range("N3").Value = 0.1113589360804189
$ws.Range("O3").Value = 0.1103713752629462
$ws.Range("P3").Value = 0.1103713752629462
$ws.Range("Q3").Value = 0.1103423393696575
$ws.Range("R3").Value = 0.1103423393696575
$ws.Range("S3").Value = 0.1103423393696575
$ws.Range("T3").Value = 0.1103423393696575
$ws.Range("U3").Value = 0.1102899996290927
$ws.Range("V3").Value = 0.1101836109344095
$ws.Range("W3").Value = 0.1101430718268682
$ws.Range("X3").Value = 0.1101430718268682
$ws.Range("Y3").Value = 0.1101430718268682
# Row 4
$ws.Range("C4").Value = 0.7541248798370361
$ws.Range("E4").Value = 5304.026819530659
$ws.Range("G4").Value = 0.1759815178231892
$ws.Range("H4").Value = 0.1704439649691365
$ws.Range("I4").Value = 0.1544716043545507
$ws.Range("J4").Value = 0.1390253213103271
$ws.Range("K4").Value = 0.1257649309364189
$ws.Range("L4").Value = 0.1163983791333658
$ws.Range("M4").Value = 0.1145779466688847
$ws.Range("N4").Value = 0.1132396693616662
$ws.Range("O4").Value = 0.1119939767195134
$ws.Range("P4").Value = 0.1119939767195134
$ws.Range("Q4").Value = 0.1114476863465585
$ws.Range("R4").Value = 0.1114476863465585
$ws.Range("S4").Value = 0.1114476863465585
$ws.Range("T4").Value = 0.1114476863465585
$ws.Range("U4").Value = 0.1114476863465585
$ws.Range("V4").Value = 0.1114476863465585
$ws.Range("W4").Value = 0.1114476863465585
$ws.Range("X4").Value = 0.1114476863465585
$ws.Range("Y4").Value = 0.1113923356633656
# Row 5
$ws.Range("C5").Value = 0.7812647819519043
$ws.Range("E5").Value = 5293.157002932115
$ws.Range("H5").Value = 0.1738735936497511
$ws.Range("I5").Value = 0.1640136354308318
$ws.Range("J5").Value = 0.1557803293636605
$ws.Range("K5").Value = 0.1557803293636605
$ws.Range("L5").Value = 0.1421768380445265
$ws.Range("M5").Value = 0.1306525455061247
$ws.Range("N5").Value = 0.1200884781812499
$ws.Range("O5").Value = 0.1142164567146343
$ws.Range("P5").Value = 0.1128810015160823
$ws.Range("Q5").Value = 0.1128810015160823
$ws.Range("R5").Value = 0.1128810015160823
$ws.Range("S5").Value = 0.1124969977363263
$ws.Range("T5").Value = 0.1120218462609039
$ws.Range("U5").Value = 0.112006882364714
$ws.Range("V5").Value = 0.1113745270724463
$ws.Range("W5").Value = 0.1113745270724463
$ws.Range("X5").Value = 0.1113199543587014
$ws.Range("Y5").Value = 0.1111804484002361
# Row 6
$ws.Range("C6").Value = 0.7500004768371582
$ws.Range("E6").Value = 5303.964419364228
$ws.Range("G6").Value = 0.1755357036184034
$ws.Range("H6").Value = 0.1640697804787035
$ws.Range("I6").Value = 0.1589793087751172
$ws.Range("J6").Value = 0.1345595289219969
$ws.Range("K6").Value = 0.1216186440732276
$ws.Range("L6").Value = 0.120186986720946
$ws.Range("M6").Value = 0.1136580202055247
$ws.Range("N6").Value = 0.1132763573958065
$ws.Range("O6").Value = 0.1129269425250347
$ws.Range("P6").Value = 0.1122000285980301
$ws.Range("Q6").Value = 0.1116684679781859
$ws.Range("R6").Value = 0.1116684679781859
$ws.Range("S6").Value = 0.1116684679781859
$ws.Range("T6").Value = 0.1115941094021108
$ws.Range("U6").Value = 0.1115527727625577
$ws.Range("V6").Value = 0.1115208954649315
$ws.Range("W6").Value = 0.1114988199640495
$ws.Range("X6").Value = 0.1113911192858524
$ws.Range("Y6").Value = 0.1113911192858524
# Row 7
$ws.Range("C7").Value = 0.7656099796295166
$ws.Range("E7").Value = 5282.017459606802
$ws.Range("G7").Value = 0.1766969141991905
$ws.Range("H7").Value = 0.164316411212265
$ws.Range("I7").Value = 0.1598154972331221
$ws.Range("J7").Value = 0.1406437875633592
$ws.Range("K7").Value = 0.1247041553975499
$ws.Range("L7").Value = 0.121775693362497
$ws.Range("M7").Value = 0.11621524269717
$ws.Range("N7").Value = 0.1159997839105884
$ws.Range("O7").Value = 0.1131130031751423
$ws.Range("P7").Value = 0.1119299582326765
$ws.Range("Q7").Value = 0.111527804651788
$ws.Range("R7").Value = 0.1111136819958433
$ws.Range("S7").Value = 0.1111136819958433
$ws.Range("T7").Value = 0.1111136819958433
$ws.Range("U7").Value = 0.1111136819958433
$ws.Range("V7").Value = 0.1110170965579736
$ws.Range("W7").Value = 0.1110170965579736
$ws.Range("X7").Value = 0.1109633033061755
$ws.Range("Y7").Value = 0.1109633033061755
# Row 8
$ws.Range("C8").Value = 0.7656407356262207
$ws.Range("E8").Value = 5220.411555874026
$ws.Range("G8").Value = 0.176567158427775
$ws.Range("H8").Value = 0.1704988410581034
$ws.Range("I8").Value = 0.1667848578834215
$ws.Range("J8").Value = 0.1600783451254889
$ws.Range("K8").Value = 0.1462654811072261
$ws.Range("L8").Value = 0.1362902297613469
$ws.Range("M8").Value = 0.1238884089341805
$ws.Range("N8").Value = 0.1202695018986898
$ws.Range("O8").Value = 0.1122873479151376
$ws.Range("P8").Value = 0.1122873479151376
$ws.Range("Q8").Value = 0.1119152366555455
$ws.Range("R8").Value = 0.1099923808320071
$ws.Range("S8").Value = 0.1099923808320071
$ws.Range("T8").Value = 0.1099923808320071
$ws.Range("U8").Value = 0.1099923808320071
$ws.Range("V8").Value = 0.1099749872789408
$ws.Range("W8").Value = 0.1098886644941043
$ws.Range("X8").Value = 0.109857368169683
$ws.Range("Y8").Value = 0.1097624084965697
# Row 9
$ws.Range("C9").Value = 0.7656099796295166
$ws.Range("E9").Value = 5346.633874085862
$ws.Range("H9").Value = 0.1681031536281684
$ws.Range("I9").Value = 0.1635900573481348
$ws.Range("J9").Value = 0.1630397854576907
$ws.Range("K9").Value = 0.1628980015650583
$ws.Range("L9").Value = 0.1495247239398385
$ws.Range("M9").Value = 0.1366424751627321
$ws.Range("N9").Value = 0.131222699331177
$ws.Range("O9").Value = 0.1212398755191479
$ws.Range("P9").Value = 0.1155117876191576
$ws.Range("Q9").Value = 0.1134719572696688
$ws.Range("R9").Value = 0.1132678214429601
$ws.Range("S9").Value = 0.1131124102144254
$ws.Range("T9").Value = 0.112862561058884
$ws.Range("U9").Value = 0.1125381242188264
$ws.Range("V9").Value = 0.1123776489882607
$ws.Range("W9").Value = 0.1123776489882607
$ws.Range("X9").Value = 0.1123364319059773
$ws.Range("Y9").Value = 0.1122228825357868
# Row 10
$ws.Range("C10").Value = 0.7499711513519287
$ws.Range("E10").Value = 5248.172324976495
$ws.Range("G10").Value = 0.1749353116568467
$ws.Range("H10").Value = 0.1743703669961806
$ws.Range("I10").Value = 0.153332217040805
$ws.Range("J10").Value = 0.1404291122959017
$ws.Range("K10").Value = 0.1223118898907158
$ws.Range("L10").Value = 0.1174028565222355
$ws.Range("M10").Value = 0.1123521752028529
$ws.Range("N10").Value = 0.1118817752150991
$ws.Range("O10").Value = 0.1118253108225645
$ws.Range("P10").Value = 0.1106301945897721
$ws.Range("Q10").Value = 0.1106301945897721
$ws.Range("R10").Value = 0.1103840750571846
$ws.Range("S10").Value = 0.1103840750571846
$ws.Range("T10").Value = 0.1103840750571846
$ws.Range("U10").Value = 0.1103840750571846
$ws.Range("V10").Value = 0.1103840750571846
$ws.Range("W10").Value = 0.1103371037618439
$ws.Range("X10").Value = 0.1103371037618439
$ws.Range("Y10").Value = 0.110303554093109
# Row 11
$ws.Range("C11").Value = 0.8281238079071045
$ws.Range("D11").Value = 3
$ws.Range("E11").Value = 6979.274756502823
$ws.Range("G11").Value = 0.1755590590951902
$ws.Range("H11").Value = 0.1668607725105126
$ws.Range("I11").Value = 0.1610084600525325
$ws.Range("J11").Value = 0.1603497614514845
$ws.Range("K11").Value = 0.1603497614514845
$ws.Range("L11").Value = 0.1603497614514845
$ws.Range("M11").Value = 0.1601369356677987
$ws.Range("N11").Value = 0.1601369356677987
$ws.Range("O11").Value = 0.1601369356677987
$ws.Range("P11").Value = 0.1601369356677987
$ws.Range("Q11").Value = 0.1601369356677987
$ws.Range("R11").Value = 0.1601369356677987
$ws.Range("S11").Value = 0.1601369356677987
$ws.Range("T11").Value = 0.1601369356677987
$ws.Range("U11").Value = 0.1601369356677987
$ws.Range("V11").Value = 0.1601369356677987
$ws.Range("W11").Value = 0.1601369356677987
$ws.Range("X11").Value = 0.1601081527615442
$ws.Range("Y11").Value = 0.1600482408675014
